$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = "LP solver (linprog or gurobi)"
$ws.Cells.Item(5, 2).Value = "gurobi"
